# Auto-update: dice_jobs_list.xlsx
# Adds a new job listing row (row 10) to the Dice jobs list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 'Backend Software Engineer'
$ws.Range("B10").Value = 'https://www.dice.com/job-detail/2f164b47-a7b0-4780-afb0-9fbee67daf87?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang'
$ws.Range("C10").Value = 'Remote or Evansville, Illinois'
$ws.Range("D10").Value = 'Contract'
$ws.Range("E10").Value = 'Market rate'
$ws.Range("F10").Value = 'Fidel Softech Ltd.'
$ws.Range("G10").Value = 'Job Description: We are looking for a highly skilled backend developer to design, build, and maintain distributed systems and APIs that power mission-critical services. This role focuses on hands-on development with TypeScript and Golang, leveraging modern frameworks and cloud-native architectures to deliver secure, scalable, and resilient solutions. Key Responsibilities Backend Development Implement APIs using TypeScript (Node.js/Apollo Server) and Golang for high-performance services. Develop GraphQL subgraphs and REST endpoints, ensuring schema integrity and efficient resolver logic. Integrate with API gateways and enforce security standards (OIDC/JWT, role-based access). System Integration Collaborate with other application teams to design and implement API integration patterns. Work with Apollo Router for federated GraphQL architectures and schema composition. Ensure robust error handling, observability, and performance optimization across services. Data Layer Build and maintain data access layers using GORM (Go) and PrismaORM/TypeORM (TypeScript). Implement data integration workflows with relational databases (PostgreSQL) and data warehouses (e.g., Snowflake). Optimize queries and ensure secure handling of sensitive data. Infrastructure & Cloud Deploy and operate services on containerized and serverless platforms (ECS/Fargate, Lambda) across AWS/Google Cloud Platform/Azure. Use IaC tools (CDK/Terraform) for environment provisioning and configuration. Implement observability practices with distributed tracing (X-Ray/OpenTelemetry), structured logging, and metrics.'
$ws.Range("H10").Value = '71.32%%'
$ws.Range("I10").Value = 'Contract Independent | On-site | Market rate'
